$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.134.14'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '1.917.07'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7928'
$ws.Range('E5').Value = '  +7.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.52'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3187'
$ws.Range('E8').Value = '  +2.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.39'
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06970'
$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08027'
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7536'
$ws.Range('E12').Value = '  -2.47%  '
$ws.Range('D13').Value = '1.902.07'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.239'
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.72'
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').Value = '30.123.46'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.05'
$ws.Range('E17').Value = '  -2.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.016'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '249.73'
$ws.Range('E19').Value = '  +3.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007827'
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').Value = '2.149.73'
$ws.Range('E22').Value = '  -0.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.993'
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '169.09'
$ws.Range('E25').Value = '  +0.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.352'
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1379'
$ws.Range('E27').Value = '  +7.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.04'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.061'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.392'
$ws.Range('E30').Value = '  +2.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.527'
$ws.Range('E31').Value = '  -1.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.375'
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.150'
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05374'
$ws.Range('E34').Value = '  +4.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.273'
$ws.Range('E35').Value = '  -2.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7416'
$ws.Range('E36').Value = '  -1.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.734'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01931'
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.794'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4479'
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.180'
$ws.Range('E41').Value = '  -2.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.86'
$ws.Range('E42').Value = '  -2.38%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.909'
$ws.Range('E43').Value = '  -4.28%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8367'
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.639'
$ws.Range('E46').Value = '  -1.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.90'
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.846'
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('D49').Value = '2.056.62'
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.57'
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '958.91'
$ws.Range('E51').Value = '  +3.26%  '
